# Added logic to reuse the current report file
# Adds two new Asset rows (ReportFileID / ReportFileURL) to the "Assets" sheet,
# mirroring the existing "GDriveMasterReportID" row (row 10), and moves the
# active-cell selection from D11 to A11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")

# Copy the formatting of the template row (row 10, columns A:D) down into the
# two new rows so the new cells pick up the same styles (s="4"/"5").
$template = $ws.Range("A10:D10")
$template.Copy()
$ws.Range("A11:D11").PasteSpecial(-4122)
$ws.Range("A12:D12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 11: ReportFileID asset
$ws.Range("A11").Value = "ReportFileID"
$ws.Range("B11").Value = "14_ReportFileID"
$ws.Range("C11").Value = "Shared"
$ws.Range("D11").Value = "Gdrive ID for the current report file"

# Row 12: ReportFileURL asset
$ws.Range("A12").Value = "ReportFileURL"
$ws.Range("B12").Value = "14_ReportFileURL"
$ws.Range("C12").Value = "Shared"
$ws.Range("D12").Value = "Gdrive URL for the current report file"

# Move the active selection from D11 to A11
$ws.Range("A11").Select() | Out-Null
